# Add a new worksheet "5-fold FS without costs " before the existing
# "cost sensitive CFS output" sheet, populate it with the 5-fold feature
# selection data (no costs), and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

$existing = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add()
$ws.Move($existing)
$ws.Name = "5-fold FS without costs "

$rows = @(
    @("CFS_ 1", "0.889274095728891", 19, "RAVLT.immediate faq_FAQBEVG faq_FAQEVENT faq_FAQFINAN faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_MMCITY mmse_objects adas_Q4SCORE adas_Q7SCORE adas_Q8SCORE TRABSCOR BNTTOTAL gd_GDMEMORY ecog_MEMORY1 ecog_PLAN5"),
    @("CFS_ 2", "0.88700240207849", 20, "RAVLT.immediate faq_FAQBEVG faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_MMFLOOR adas_Q4SCORE adas_Q1SCORE adas_Q10SCORE TRABSCOR BNTTOTAL AVDELTOT gd_GDMEMORY ecog_MEMORY1"),
    @("CFS_ 3", "0.880376556298963", 24, "faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_MMDATE adas_Q4SCORE adas_Q1SCORE adas_Q7SCORE CATANIMSC TRABSCOR BNTTOTAL gd_GDDROP gd_GDMEMORY gd_GDWORTH ecog_LANG7 ecog_MEMORY1 ecog_MEMORY2 ecog_MEMORY4 ecog_ORGAN6"),
    @("CFS_ 4", "0.867390161612428", 16, "faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_objects adas_Q4SCORE adas_Q7SCORE CATANIMSC BNTTOTAL AVDELTOT ecog_LANG9 ecog_PLAN3"),
    @("CFS_ 5", "0.884976761667935", 18, "RAVLT.immediate faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL mmse_MMDAY mmse_objects adas_Q4SCORE adas_Q1SCORE adas_Q8SCORE CATANIMSC TRABSCOR gd_GDMEMORY ecog_PLAN3"),
    @("Consistency_ 1", "0.8702725037447", 16, "RAVLT.immediate faq_FAQBEVG faq_FAQEVENT faq_FAQFINAN faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_MMCITY mmse_objects adas_Q4SCORE adas_Q7SCORE adas_Q8SCORE TRABSCOR BNTTOTAL gd_GDMEMORY ecog_MEMORY1 ecog_PLAN5"),
    @("Consistency_ 2", "0.864447060409649", 16, "RAVLT.immediate faq_FAQBEVG faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_MMFLOOR adas_Q4SCORE adas_Q1SCORE adas_Q10SCORE TRABSCOR BNTTOTAL AVDELTOT gd_GDMEMORY ecog_MEMORY1"),
    @("Consistency_ 3", "0.871944448586633", 17, "faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_MMDATE adas_Q4SCORE adas_Q1SCORE adas_Q7SCORE CATANIMSC TRABSCOR BNTTOTAL gd_GDDROP gd_GDMEMORY gd_GDWORTH ecog_LANG7 ecog_MEMORY1 ecog_MEMORY2 ecog_MEMORY4 ecog_ORGAN6"),
    @("Consistency_ 4", "0.863316617842934", 19, "faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_objects adas_Q4SCORE adas_Q7SCORE CATANIMSC BNTTOTAL AVDELTOT ecog_LANG9 ecog_PLAN3"),
    @("Consistency_ 5", "0.86742070423488", 13, "RAVLT.immediate faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL mmse_MMDAY mmse_objects adas_Q4SCORE adas_Q1SCORE adas_Q8SCORE CATANIMSC TRABSCOR gd_GDMEMORY ecog_PLAN3"),
    @("Boruta_ 1", "0.891311709377254", 20, "RAVLT.immediate moca_language faq_FAQBEVG faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_objects adas_Q4SCORE adas_Q1SCORE adas_Q7SCORE gd_GDMEMORY ecog_DIVATT4 ecog_LANG1 ecog_MEMORY1"),
    @("Boruta_ 2", "0.884174293292747", 19, "RAVLT.immediate faq_FAQBEVG faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV adas_Q4SCORE adas_Q1SCORE adas_Q10SCORE CATANIMSC AVDELTOT gd_GDMEMORY ecog_LANG9 ecog_MEMORY1"),
    @("Boruta_ 3", "0.885283906144912", 23, "faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_MMDATE adas_Q4SCORE adas_Q1SCORE CATANIMSC TRABSCOR BNTTOTAL gd_GDMEMORY ecog_DIVATT4 ecog_LANG7 ecog_MEMORY1 ecog_MEMORY2 ecog_MEMORY3 ecog_MEMORY4 ecog_MEMORY6"),
    @("Boruta_ 4", "0.877143304749275", 20, "RAVLT.immediate faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL faq_FAQTV mmse_objects adas_Q4SCORE adas_Q1SCORE adas_Q7SCORE CATANIMSC BNTTOTAL AVDELTOT ecog_LANG9 ecog_PLAN3 ecog_PLAN4"),
    @("Boruta_ 5", "0.884387143544527", 20, "RAVLT.immediate faq_FAQEVENT faq_FAQFINAN faq_FAQFORM faq_FAQGAME faq_FAQMEAL faq_FAQREM faq_FAQSHOP faq_FAQTRAVL adas_Q4SCORE adas_Q1SCORE adas_Q7SCORE adas_Q8SCORE TRABSCOR gd_GDMEMORY gd_GDWORTH ecog_LANG7 ecog_LANG9 ecog_MEMORY5 ecog_PLAN3")
)

$lastRow = $rows.Length + 1

# The AUC column (B) holds numeric-looking values that were entered/stored
# as text in the source workbook. Force text formatting while writing them,
# then restore the default "Normal" style so no stray number-format style
# is left attached to the cells.
$aucRange = $ws.Range("B2:B$lastRow")
$aucRange.NumberFormat = "@"

# Reproduce the original authoring order: the B/C/D (AUC, feature count,
# feature list) columns -- header first, then row by row -- were filled in
# before the A (Algorithm and Fold) column.
$ws.Cells.Item(1, 2).Value = "AUC"
$ws.Cells.Item(1, 3).Value = "Number of features"
$ws.Cells.Item(1, 4).Value = "List of Features"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $ws.Cells.Item($r + 2, 2).Value = $row[1]
    $ws.Cells.Item($r + 2, 3).Value = $row[2]
    $ws.Cells.Item($r + 2, 4).Value = $row[3]
}

$ws.Cells.Item(1, 1).Value = "Algorithm and Fold"
for ($r = 0; $r -lt $rows.Length; $r++) {
    $ws.Cells.Item($r + 2, 1).Value = $rows[$r][0]
}

$aucRange.Style = "Normal"

$ws.Activate()
$ws.Range("F15").Select()
